$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($cell, [string]$text)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextCell $ws.Range("D2") '98.370.15'
Set-TextCell $ws.Range("E2") '  +0.72%  '
Set-TextCell $ws.Range("D3") '3.325.91'
Set-TextCell $ws.Range("E3") '  +6.09%  '
Set-TextCell $ws.Range("E4") '  +0.08%  '
Set-TextCell $ws.Range("D5") '257.56'
Set-TextCell $ws.Range("E5") '  +6.25%  '
Set-TextCell $ws.Range("D6") '625.11'
Set-TextCell $ws.Range("E6") '  +2.54%  '
Set-TextCell $ws.Range("E7") '  +24.49%  '
Set-TextCell $ws.Range("D8") '0.387'
Set-TextCell $ws.Range("E8") '  +0.87%  '
Set-TextCell $ws.Range("E9") '  +0.02%  '
Set-TextCell $ws.Range("D10") '0.861'
Set-TextCell $ws.Range("E10") '  +9.99%  '
Set-TextCell $ws.Range("D11") '3.321.31'
Set-TextCell $ws.Range("E11") '  +5.96%  '
Set-TextCell $ws.Range("D12") '0.198'
Set-TextCell $ws.Range("E12") '  +0.17%  '
Set-TextCell $ws.Range("D13") '36.91'
Set-TextCell $ws.Range("E13") '  +9.31%  '
Set-TextCell $ws.Range("D14") '98.191.33'
Set-TextCell $ws.Range("E14") '  +0.90%  '
Set-TextCell $ws.Range("D15") '0.0000247'
Set-TextCell $ws.Range("E15") '  +3.20%  '
Set-TextCell $ws.Range("D16") '3.938.89'
Set-TextCell $ws.Range("E16") '  +5.96%  '
Set-TextCell $ws.Range("D17") '5.49'
Set-TextCell $ws.Range("E17") '  +0.97%  '
Set-TextCell $ws.Range("D18") '3.331.71'
Set-TextCell $ws.Range("E18") '  +6.33%  '
Set-TextCell $ws.Range("D19") '3.53'
Set-TextCell $ws.Range("E19") '  +3.28%  '
Set-TextCell $ws.Range("D20") '15.07'
Set-TextCell $ws.Range("E20") '  +4.80%  '
Set-TextCell $ws.Range("D21") '486.99'
Set-TextCell $ws.Range("E21") '  -6.08%  '
Set-TextCell $ws.Range("D22") '6.04'
Set-TextCell $ws.Range("E22") '  +7.34%  '
Set-TextCell $ws.Range("D23") '0.0000210'
Set-TextCell $ws.Range("E23") '  +9.96%  '
Set-TextCell $ws.Range("D24") '9.32'
Set-TextCell $ws.Range("E24") '  +7.25%  '
Set-TextCell $ws.Range("D25") '5.60'
Set-TextCell $ws.Range("E25") '  +3.11%  '
Set-TextCell $ws.Range("D26") '88.55'
Set-TextCell $ws.Range("E26") '  +0.35%  '
Set-TextCell $ws.Range("D27") '11.86'
Set-TextCell $ws.Range("E27") '  +3.05%  '
Set-TextCell $ws.Range("D28") '3.505.84'
Set-TextCell $ws.Range("E28") '  +6.53%  '
Set-TextCell $ws.Range("E29") '  +18.39%  '
Set-TextCell $ws.Range("E30") '  -0.05%  '
Set-TextCell $ws.Range("E31") '  +10.53%  '
Set-TextCell $ws.Range("E32") '  +12.35%  '
Set-TextCell $ws.Range("D33") '0.999'
Set-TextCell $ws.Range("E33") '  -0.04%  '
Set-TextCell $ws.Range("D34") '9.58'
Set-TextCell $ws.Range("E34") '  +7.84%  '
Set-TextCell $ws.Range("D35") '27.70'
Set-TextCell $ws.Range("E35") '  +4.23%  '
Set-TextCell $ws.Range("D36") '0.149'
Set-TextCell $ws.Range("E36") '  -0.70%  '
Set-TextCell $ws.Range("D37") '7.27'
Set-TextCell $ws.Range("E37") '  +0.50%  '
Set-TextCell $ws.Range("D38") '1.94'
Set-TextCell $ws.Range("E38") '  +3.78%  '
Set-TextCell $ws.Range("D39") '494.13'
Set-TextCell $ws.Range("E39") '  +5.47%  '
Set-TextCell $ws.Range("B40") 'WhiteBITCoin'
Set-TextCell $ws.Range("C40") 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
Set-TextCell $ws.Range("D40") '24.84'
Set-TextCell $ws.Range("E40") '  +2.11%  '
Set-TextCell $ws.Range("B41") 'PolygonEcosystemToken'
Set-TextCell $ws.Range("C41") 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
Set-TextCell $ws.Range("D41") '0.457'
Set-TextCell $ws.Range("E41") '  +5.62%  '
Set-TextCell $ws.Range("E42") '  +4.88%  '
Set-TextCell $ws.Range("D43") '1.25'
Set-TextCell $ws.Range("E43") '  +3.72%  '
Set-TextCell $ws.Range("D44") '3.28'
Set-TextCell $ws.Range("E44") '  +6.33%  '
Set-TextCell $ws.Range("E45") '  -0.02%  '
Set-TextCell $ws.Range("D46") '0.777'
Set-TextCell $ws.Range("E46") '  +12.12%  '
Set-TextCell $ws.Range("D47") '159.08'
Set-TextCell $ws.Range("E47") '  -2.05%  '
Set-TextCell $ws.Range("D48") '1.92'
Set-TextCell $ws.Range("E48") '  +1.33%  '
Set-TextCell $ws.Range("D49") '0.841'
Set-TextCell $ws.Range("E49") '  +8.97%  '
Set-TextCell $ws.Range("D50") '4.59'
Set-TextCell $ws.Range("E50") '  +2.99%  '
Set-TextCell $ws.Range("D51") '45.48'
Set-TextCell $ws.Range("E51") '  +3.09%  '
